$wb = $excel.ActiveWorkbook

# --- Sheet "Hoja1": update the conversion text in A1 ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$newText = "Conversión del día 💰`r`n✅ Dólar paralelo: 68`r`n`r`nBinance`r`n✅ 1000 Bs = 8.37 = 34309.62 pesos`r`n✅ 34309.62 pesos = 8.34 = 980.04 Bs`r`n`r`nPromedio competencia`r`n✅ Tasa pesos: 20`r`n✅ Tasa Bs: 20`r`n✅ % Ganancia: 20%"
$wsHoja1.Range("A1").Value = $newText

# --- Sheet "tasas": update N10, O10, N12 ---
$wsTasas = $wb.Worksheets.Item("tasas")
$wsTasas.Range("N10").Value = 119.5
$wsTasas.Range("O10").Value = 4100
$wsTasas.Range("N12").Value = 4113.5
